$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 521
$ws.Range("I5").Value = 401.5
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 401.5
$ws.Range("L5").Value = 999
$ws.Range("M5").Value = -286.5
$ws.Range("N5").Value = -1229
$ws.Range("H6").Value = 165.16667
$ws.Range("I6").Value = 165.16667
$ws.Range("K6").Value = 495.50001
$ws.Range("M6").Value = -383.50001
$ws.Range("H19").Value = 1401.1666
$ws.Range("J19").Value = 1458.75
$ws.Range("L19").Value = 1458.75
$ws.Range("N19").Value = -1808.75
$ws.Range("H86").Value = 1499
$ws.Range("I86").Value = 1499
$ws.Range("K86").Value = 1499
$ws.Range("M86").Value = -376
$ws.Range("H89").Value = 1499
$ws.Range("I89").Value = 1499
$ws.Range("K89").Value = 7495
$ws.Range("M89").Value = -1879
$ws.Range("H138").Value = 2134.2
$ws.Range("I138").Value = 2134.2
$ws.Range("K138").Value = 6402.599999999999
$ws.Range("M138").Value = -1262.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1003.3333
$ws.Range("I45").Value = 1003.3333
$ws.Range("K45").Value = 1003.3333
$ws.Range("M45").Value = -626.3333
$ws.Range("H61").Value = 3044.6924
$ws.Range("I61").Value = 3044.6924
$ws.Range("K61").Value = 3044.6924
$ws.Range("M61").Value = -2832.6924
$ws.Range("H74").Value = 3497.4
$ws.Range("I74").Value = 4500
$ws.Range("K74").Value = 4500
$ws.Range("M74").Value = -3626
$ws.Range("H77").Value = 3497.4
$ws.Range("I77").Value = 4500
$ws.Range("K77").Value = 22500
$ws.Range("M77").Value = -18132
$ws.Range("H110").Value = 11477.5
$ws.Range("J110").Value = 1949.5
$ws.Range("L110").Value = 1949.5
$ws.Range("N110").Value = -6039.5
$ws.Range("H122").Value = 1422.3334
$ws.Range("I122").Value = 1422.3334
$ws.Range("K122").Value = 4267.0002
$ws.Range("M122").Value = -1817.0002
$ws.Range("H132").Value = 1546.2903
$ws.Range("I132").Value = 1505.64
$ws.Range("J132").Value = 1715.6666
$ws.Range("K132").Value = 4516.92
$ws.Range("L132").Value = 5146.9998
$ws.Range("M132").Value = -1986.92
$ws.Range("N132").Value = -10206.9998
$ws.Range("H136").Value = 3044.6924
$ws.Range("I136").Value = 3044.6924
$ws.Range("K136").Value = 9134.0772
$ws.Range("M136").Value = -6584.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5298.6
$ws.Range("I86").Value = 5298.6
$ws.Range("K86").Value = 5298.6
$ws.Range("M86").Value = -4175.6
$ws.Range("H89").Value = 5298.6
$ws.Range("I89").Value = 5298.6
$ws.Range("K89").Value = 26493
$ws.Range("M89").Value = -20877
$ws.Range("H105").Value = 2815.2856
$ws.Range("I105").Value = 1598.6
$ws.Range("K105").Value = 1598.6
$ws.Range("M105").Value = 148.4000000000001
$ws.Range("H107").Value = 1658.8572
$ws.Range("I107").Value = 1656.8
$ws.Range("K107").Value = 1656.8
$ws.Range("M107").Value = 263.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 254
$ws.Range("I2").Value = 254
$ws.Range("K2").Value = 254
$ws.Range("M2").Value = -141
$ws.Range("H5").Value = 419.6
$ws.Range("I5").Value = 419.6
$ws.Range("K5").Value = 419.6
$ws.Range("M5").Value = -307.6
$ws.Range("H16").Value = 71429540
$ws.Range("I16").Value = 100000790
$ws.Range("J16").Value = 1399.5
$ws.Range("K16").Value = 100000790
$ws.Range("L16").Value = 1399.5
$ws.Range("M16").Value = -100000503
$ws.Range("N16").Value = -1973.5
$ws.Range("H31").Value = 4676.8335
$ws.Range("I31").Value = 3569.4443
$ws.Range("K31").Value = 3569.4443
$ws.Range("M31").Value = -3274.4443
$ws.Range("H34").Value = 4676.8335
$ws.Range("I34").Value = 3569.4443
$ws.Range("K34").Value = 3569.4443
$ws.Range("M34").Value = -3367.4443
$ws.Range("H105").Value = 849.5
$ws.Range("I105").Value = 849
$ws.Range("K105").Value = 849
$ws.Range("M105").Value = 898
$ws.Range("H113").Value = 71429540
$ws.Range("I113").Value = 100000790
$ws.Range("J113").Value = 1399.5
$ws.Range("K113").Value = 100000790
$ws.Range("L113").Value = 1399.5
$ws.Range("M113").Value = -99998620
$ws.Range("N113").Value = -5739.5
$ws.Range("H134").Value = 3093.25
$ws.Range("I134").Value = 3093.25
$ws.Range("K134").Value = 9279.75
$ws.Range("M134").Value = -6744.75
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 413
$ws.Range("I10").Value = 217.66667
$ws.Range("K10").Value = 653.00001
$ws.Range("M10").Value = -514.00001
$ws.Range("H13").Value = 322
$ws.Range("I13").Value = 45.5
$ws.Range("J13").Value = 875
$ws.Range("K13").Value = 136.5
$ws.Range("L13").Value = 2625
$ws.Range("M13").Value = 31.5
$ws.Range("N13").Value = -2961
$ws.Range("H21").Value = 2000.5
$ws.Range("I21").Value = 2000.5
$ws.Range("K21").Value = 6001.5
$ws.Range("M21").Value = -5828.5
$ws.Range("H25").Value = 949.8
$ws.Range("I25").Value = 687.25
$ws.Range("K25").Value = 2061.75
$ws.Range("M25").Value = -1892.75
$ws.Range("H30").Value = 949.8
$ws.Range("I30").Value = 687.25
$ws.Range("K30").Value = 2061.75
$ws.Range("M30").Value = -1959.75
$ws.Range("H32").Value = 999
$ws.Range("J32").Value = 999
$ws.Range("L32").Value = 2997
$ws.Range("N32").Value = -3563
$ws.Range("H70").Value = 17142.857
$ws.Range("I70").Value = 13300
$ws.Range("J70").Value = 17783.334
$ws.Range("K70").Value = 39900
$ws.Range("L70").Value = 53350.00199999999
$ws.Range("M70").Value = -39585
$ws.Range("N70").Value = -53980.00199999999
$ws.Range("H73").Value = 17142.857
$ws.Range("I73").Value = 13300
$ws.Range("J73").Value = 17783.334
$ws.Range("K73").Value = 39900
$ws.Range("L73").Value = 53350.00199999999
$ws.Range("M73").Value = -38808
$ws.Range("N73").Value = -55534.00199999999
$ws.Range("H87").Value = 1633
$ws.Range("I87").Value = 1633
$ws.Range("K87").Value = 4899
$ws.Range("M87").Value = -3651
$ws.Range("H90").Value = 1633
$ws.Range("I90").Value = 1633
$ws.Range("K90").Value = 14697
$ws.Range("M90").Value = -8457
$ws.Range("H97").Value = 1935.4
$ws.Range("I97").Value = 338.75
$ws.Range("J97").Value = 2999.8333
$ws.Range("K97").Value = 1016.25
$ws.Range("L97").Value = 8999.499899999999
$ws.Range("M97").Value = -520.25
$ws.Range("N97").Value = -9991.499899999999
$ws.Range("H114").Value = 1086.7142
$ws.Range("I114").Value = 1283.4
$ws.Range("J114").Value = 595
$ws.Range("K114").Value = 3850.2
$ws.Range("L114").Value = 1785
$ws.Range("M114").Value = -596.2000000000003
$ws.Range("N114").Value = -8293
$ws.Range("H122").Value = 451.66666
$ws.Range("J122").Value = 502.5
$ws.Range("L122").Value = 4522.5
$ws.Range("N122").Value = -9422.5
$ws.Range("H131").Value = 2120.375
$ws.Range("I131").Value = 1893.75
$ws.Range("K131").Value = 5681.25
$ws.Range("M131").Value = -641.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11787070
$ws.Range("I11").Value = 2336499.2
$ws.Range("J11").Value = 18874998
$ws.Range("K11").Value = 2336499.2
$ws.Range("L11").Value = 18874998
$ws.Range("M11").Value = -2336360.2
$ws.Range("N11").Value = -18875276
$ws.Range("H113").Value = 1862.4445
$ws.Range("I113").Value = 1608.8572
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 1608.8572
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = 561.1428000000001
$ws.Range("N113").Value = -7090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 167
$ws.Range("N19").Value = ""
$ws.Range("H25").Value = 2500
$ws.Range("I25").Value = 2500
$ws.Range("K25").Value = 2500
$ws.Range("M25").Value = -2270
$ws.Range("H31").Value = 2162
$ws.Range("I31").Value = 1452.5
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1452.5
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1204.5
$ws.Range("N31").Value = -5496
$ws.Range("H46").Value = 833.6
$ws.Range("I46").Value = 833.6
$ws.Range("K46").Value = 833.6
$ws.Range("M46").Value = -645.6
$ws.Range("H82").Value = 4782.091
$ws.Range("J82").Value = 6780.8
$ws.Range("L82").Value = 6780.8
$ws.Range("N82").Value = -7502.8
$ws.Range("H85").Value = 4782.091
$ws.Range("J85").Value = 6780.8
$ws.Range("L85").Value = 6780.8
$ws.Range("N85").Value = -9276.799999999999
$ws.Range("H132").Value = 3502.75
$ws.Range("I132").Value = 3502.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10508.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7978.25
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25032248
$ws.Range("J2").Value = 40000
$ws.Range("L2").Value = 40000
$ws.Range("N2").Value = -40224
